$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.353686666666667
$ws.Range("H2").Value = 4.06106
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.502378
$ws.Range("N2").Value = 1.507134
$ws.Range("O2").Value = 0.0189742916423209
$ws.Range("P2").Value = 0.0189742916423209
$ws.Range("Q2").Value = 0.6800624002266666
$ws.Range("R2").Value = 6.12056160204
$ws.Range("S2").Value = 0.0189742916423209
$ws.Range("T2").Value = 0.0189742916423209

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.353686666666667
$ws.Range("H3").Value = 4.06106
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6023626666666666
$ws.Range("N3").Value = 1.807088
$ws.Range("O3").Value = 0.02275060793223323
$ws.Range("P3").Value = 0.02275060793223323
$ws.Range("Q3").Value = 0.8154103103644443
$ws.Range("R3").Value = 7.33869279328
$ws.Range("S3").Value = 0.02275060793223323
$ws.Range("T3").Value = 0.02275060793223323

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.353686666666667
$ws.Range("H4").Value = 4.06106
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.474800666666667
$ws.Range("N4").Value = 16.424402
$ws.Range("O4").Value = 0.2067774952981744
$ws.Range("P4").Value = 0.2067774952981744
$ws.Range("Q4").Value = 7.411164665124445
$ws.Range("R4").Value = 66.70048198612001
$ws.Range("S4").Value = 0.2067774952981744
$ws.Range("T4").Value = 0.2067774952981744

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.353686666666667
$ws.Range("H5").Value = 4.06106
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.012936333333333
$ws.Range("N5").Value = 15.038809
$ws.Range("O5").Value = 0.1893333624741797
$ws.Range("P5").Value = 0.1893333624741797
$ws.Range("Q5").Value = 6.785945075282222
$ws.Range("R5").Value = 61.07350567754001
$ws.Range("S5").Value = 0.1893333624741797
$ws.Range("T5").Value = 0.1893333624741797

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.353686666666667
$ws.Range("H6").Value = 4.06106
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.423250333333333
$ws.Range("N6").Value = 22.269751
$ws.Range("O6").Value = 0.2803684014001858
$ws.Range("P6").Value = 0.2803684014001858
$ws.Range("Q6").Value = 10.04875499956222
$ws.Range("R6").Value = 90.43879499606001
$ws.Range("S6").Value = 0.2803684014001858
$ws.Range("T6").Value = 0.2803684014001858

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.353686666666667
$ws.Range("H7").Value = 4.06106
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.461044333333334
$ws.Range("N7").Value = 22.383133
$ws.Range("O7").Value = 0.281795841252906
$ws.Range("P7").Value = 0.281795841252906
$ws.Range("Q7").Value = 10.09991623344222
$ws.Range("R7").Value = 90.89924610098001
$ws.Range("S7").Value = 0.281795841252906
$ws.Range("T7").Value = 0.281795841252906
